# Apply updated "想去人数" (column F) figures across the relevant sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 9118
$ws1.Range("F8").Value  = 149
$ws1.Range("F9").Value  = 214
$ws1.Range("F11").Value = 381
$ws1.Range("F14").Value = 412
$ws1.Range("F15").Value = 11660
$ws1.Range("F16").Value = 11660
$ws1.Range("F37").Value = 959
$ws1.Range("F38").Value = 4163
$ws1.Range("F39").Value = 319
$ws1.Range("F40").Value = 3048
$ws1.Range("F41").Value = 1288
$ws1.Range("F45").Value = 444
$ws1.Range("F47").Value = 119
$ws1.Range("F48").Value = 180
$ws1.Range("F50").Value = 114

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 18
$ws2.Range("F18").Value = 5

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 18
$ws4.Range("F9").Value  = 9118
$ws4.Range("F13").Value = 214
$ws4.Range("F15").Value = 381
$ws4.Range("F18").Value = 11660
$ws4.Range("F38").Value = 959
$ws4.Range("F39").Value = 5
$ws4.Range("F40").Value = 4163
$ws4.Range("F41").Value = 319
$ws4.Range("F42").Value = 3048
$ws4.Range("F43").Value = 1288
$ws4.Range("F47").Value = 119
$ws4.Range("F48").Value = 180
$ws4.Range("F50").Value = 114
